$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.960.54'
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").Value = '3.853.73'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''598.66'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '''168.44'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").Value = '3.857.25'
$ws.Range("E7").Value = '  -1.76%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '''0.528'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("E10").Value = '  -4.61%  '
$ws.Range("D11").Value = '''6.42'
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("D12").Value = '''0.455'
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("D13").Value = '''0.0000258'
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("D14").Value = '''36.82'
$ws.Range("E14").Value = '  -2.20%  '
$ws.Range("D15").Value = '4.498.39'
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").Value = '3.850.20'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '68.041.15'
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("D18").Value = '''17.98'
$ws.Range("E18").Value = '  +3.65%  '
$ws.Range("D19").Value = '''7.32'
$ws.Range("E19").Value = '  -2.22%  '
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("D22").Value = '''465.18'
$ws.Range("E22").Value = '  -6.39%  '
$ws.Range("D23").Value = '''0.734'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = '''0.0000159'
$ws.Range("E24").Value = '  -4.56%  '
$ws.Range("D25").Value = '''82.81'
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("D27").Value = '''12.00'
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").Value = '''9.93'
$ws.Range("E29").Value = '  -2.89%  '
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").Value = '4.001.50'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").Value = '''7.66'
$ws.Range("E32").Value = '  -1.43%  '
$ws.Range("E33").Value = '  -3.59%  '
$ws.Range("D34").Value = '''31.07'
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("D35").Value = '''9.45'
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("D36").Value = '3.814.24'
$ws.Range("E36").Value = '  -1.98%  '
$ws.Range("E37").Value = '  -2.86%  '
$ws.Range("D38").Value = '''3.63'
$ws.Range("E38").Value = '  +10.93%  '
$ws.Range("E39").Value = '  -2.12%  '
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").Value = '''5.89'
$ws.Range("E41").Value = '  -1.99%  '
$ws.Range("D42").Value = '''0.998'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").Value = '''0.311'
$ws.Range("E43").Value = '  -3.34%  '
$ws.Range("D44").Value = '''422.55'
$ws.Range("E44").Value = '  -1.97%  '
$ws.Range("D45").Value = '''1.97'
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("D46").Value = '''0.000294'
$ws.Range("E46").Value = '  +4.99%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '''8.60'
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '''47.09'
$ws.Range("E49").Value = '  -2.10%  '
$ws.Range("D50").Value = '''26.46'
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("D51").Value = '''142.29'
$ws.Range("E51").Value = '  -0.48%  '

Write-Output "Applied cryptos update."
